# break out stock.yaml completed
# - convert bsecode (column D) on existing rows 213-221 from text to numeric
# - append new rows 222-227 (fresh "day" snapshot rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- Fix up bsecode column (D) for rows 213-221: was stored as text, should be numeric ---
$bsecodes = @{
    213 = 532500
    214 = 500790
    215 = 532868
    216 = 500440
    217 = 540611
    218 = 540222
    219 = 500104
    220 = 500547
    221 = 500049
}
foreach ($row in $bsecodes.Keys) {
    $ws.Range("D$row").Value = $bsecodes[$row]
}

# --- Append new rows 222-227 ---
$newRows = @(
    @{ row=222; sr=1; nsecode="NAUKRI";     name="Info Edge (india) Limited";                     bsecode="532777"; per_chg=2.48;  close=7187.4;  volume=360587;    timeframe="day"; dt="26/07/2024 11:39:15" },
    @{ row=223; sr=2; nsecode="HAVELLS";    name="Havells India Limited";                          bsecode="517354"; per_chg=0.65;  close=1838.65; volume=814897;    timeframe="day"; dt="26/07/2024 11:39:15" },
    @{ row=224; sr=3; nsecode="ICICIBANK";  name="Icici Bank Limited";                              bsecode="532174"; per_chg=0.78;  close=1207.2;  volume=19309265;  timeframe="day"; dt="26/07/2024 11:39:15" },
    @{ row=225; sr=4; nsecode="SBICARD";    name="SBI Cards & Payment Services Ltd";                bsecode="543066"; per_chg=-1.2;  close=721.7;   volume=3001601;   timeframe="day"; dt="26/07/2024 11:39:15" },
    @{ row=226; sr=5; nsecode="POWERGRID";  name="Power Grid Corporation Of India Limited";         bsecode="532898"; per_chg=1.4;   close=344.2;   volume=11468068;  timeframe="day"; dt="26/07/2024 11:39:15" },
    @{ row=227; sr=6; nsecode="M&MFIN";     name="Mahindra & Mahindra Financial Services Limited";  bsecode="532720"; per_chg=1.31;  close=293.65;  volume=2051866;   timeframe="day"; dt="26/07/2024 11:39:15" }
)

foreach ($r in $newRows) {
    $row = $r.row
    $ws.Range("A$row").Value = $r.sr
    $ws.Range("B$row").Value = $r.nsecode
    $ws.Range("C$row").Value = $r.name

    # bsecode must stay a text value (leading/trailing-zero stock codes), not get
    # auto-coerced to a number -> force the text number format before assignment,
    # then strip the custom style back off so the cell carries no style index.
    $ws.Range("D$row").NumberFormat = "@"
    $ws.Range("D$row").Value = $r.bsecode
    $ws.Range("D$row").Style = "Normal"

    $ws.Range("E$row").Value = $r.per_chg
    $ws.Range("F$row").Value = $r.close
    $ws.Range("G$row").Value = $r.volume
    $ws.Range("H$row").Value = $r.timeframe
    $ws.Range("I$row").Value = $r.dt
}
